$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the sequence-diagram exception step text: "passo 5" -> "passo 4"
$ws.Range("B14").Value = "Exceção      [Cancela edição]     (passo 4)"

# Match the author's final selection (merged range B14:B18)
$ws.Range("B14:B18").Select() | Out-Null
